# Applies the cryptos-list refresh described by the commit message:
#   "Updated cryptos list on Tue Oct 22 21:30:44 UTC 2024 with GitHub Actions"
#
# The sheet stores every data cell (Coin/Link/Price/Volume) as text, even
# when the "Price" column looks numeric (e.g. "595.10", "0.362"). Some of
# those price strings (e.g. "67.416.27", "3.113.37") use "." as a thousands
# separator and are not valid numbers at all. Plain `Range.Value = "595.10"`
# would let Excel auto-coerce the numeric-looking strings into real numbers
# (dropping trailing zeros / reformatting), so those assignments are given a
# leading apostrophe - exactly like a user typing '595.10 into the Excel UI -
# to force a literal text entry and preserve the original formatting.
#
# Row 41/42 additionally swap which coin occupies which row (WhiteBITCoin
# moves to row 41, Stacks moves to row 42), so Coin/Link/Price/Volume are all
# rewritten for those two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.416.27"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "2.633.15"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'595.10"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "'167.42"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.60%  "
$ws.Range("D9").Value = "2.632.56"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").Value = "'0.362"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'5.23"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "'27.64"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "3.113.37"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "'0.0000181"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "67.446.26"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "2.628.05"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").Value = "'356.88"
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").Value = "'4.67"
$ws.Range("E23").Value = "  -3.62%  "
$ws.Range("D24").Value = "'1.94"
$ws.Range("E24").Value = "  -5.18%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "'10.27"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "'69.85"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("D31").Value = "'547.81"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").Value = "'7.93"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("E33").Value = "  -3.75%  "
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("E35").Value = "  +4.10%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -4.37%  "
$ws.Range("D38").Value = "'157.70"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").Value = "'18.99"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "'18.30"
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.82"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  -4.47%  "
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("D47").Value = "'152.89"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "'0.579"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("E51").Value = "  -1.20%  "
